$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Formula Samples" sheet: populate it with the same kind of sample data
#    as "Constant Samples", but using formulas instead of literal constants.
# ---------------------------------------------------------------------------
$wsFormula = $wb.Worksheets.Item("Formula Samples")
$wsFormula.Activate()

$wsFormula.Columns.Item(1).ColumnWidth = 13.8
$wsFormula.Columns.Item(2).ColumnWidth = 14.8

$wsFormula.Range("A1").Value = "Sample Data"
$wsFormula.Range("B1").Value = "Sample Value"
$wsFormula.Range("E1").Value = "Notes"

$wsFormula.Range("A2").Value = "Year"
$wsFormula.Range("B2").Value = 2020
$wsFormula.Range("E2").Value = "Remember that constant comparison will only compare the formula, not the computed values."

$wsFormula.Range("A3").Value = "Decade"
$wsFormula.Range("B3").Formula = "=B2 / 10"

$wsFormula.Range("A4").Value = "Century"
$wsFormula.Range("B4").Formula = "=  ROUNDUP(B2 / 100, 0)"

$wsFormula.Range("A6").Value = "Meter"
$wsFormula.Range("B6").Value = 1000

$wsFormula.Range("A7").Value = "Kilometer"
$wsFormula.Range("B7").Formula = "= B6 / 1000"

$wsFormula.Range("A8").Value = "Feet"
$wsFormula.Range("B8").Formula = "=B6 * 3.28084"

$wsFormula.Range("A9").Value = "Inch"
$wsFormula.Range("B9").Formula = "=B6 * 39.3701"
$wsFormula.Range("E9").Value = "The Inch (B9) has alt_cells in the rubric, which allows the submission to be compared to other cells."

$wsFormula.Range("C10").Formula = "=B6 * 39.37"

# Rubric comments on the formula cells (mirrors the comments already present
# on "Constant Samples", adapted for type: formula).
$wsFormula.Range("B3").AddComment("rubric:`n score: 1.5`n type: formula`n`n") | Out-Null
$wsFormula.Range("B4").AddComment("rubric:`n score: 1.5`n type: formula`n`n") | Out-Null
$wsFormula.Range("B7").AddComment("rubric:`n score: 1`n type: formula`n`n") | Out-Null
$wsFormula.Range("B8").AddComment("rubric:`n score: 1.5`n type: formula`n") | Out-Null
$wsFormula.Range("B9").AddComment("rubric:`n score: 1.5`n type: formula`nalt_cells:`n - C10`n") | Out-Null

$excel.ActiveWindow.Zoom = 159
$wsFormula.Range("B9").Select()

# ---------------------------------------------------------------------------
# 2) "Constant Samples": the active-cell marker moved off this sheet, so it
#    is no longer the selected tab; just move the remembered selection.
# ---------------------------------------------------------------------------
$wsConstant = $wb.Worksheets.Item("Constant Samples")
$wsConstant.Activate()
$wsConstant.Range("E9").Select()

# ---------------------------------------------------------------------------
# 3) "Constant Samples_CheckOrder": select the whole data range and switch
#    the page setup to portrait orientation.
# ---------------------------------------------------------------------------
$wsConstantCheck = $wb.Worksheets.Item("Constant Samples_CheckOrder")
$wsConstantCheck.Activate()
$wsConstantCheck.Range("A1:B8").Select()
$wsConstantCheck.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4) "Formula Samples_CheckOrder": populate with the matching cell-order
#    list (mirrors "Constant Samples_CheckOrder") and make it the final
#    active sheet/selection.
# ---------------------------------------------------------------------------
$wsFormulaCheck = $wb.Worksheets.Item("Formula Samples_CheckOrder")
$wsFormulaCheck.Activate()

$wsFormulaCheck.Range("A2").Value = 1
$wsFormulaCheck.Range("B2").Value = "B2"
$wsFormulaCheck.Range("C2").Value = "> This cell has no rubric, so it won't be processed."

$wsFormulaCheck.Range("A3").Value = 2
$wsFormulaCheck.Range("B3").Value = "B3"

$wsFormulaCheck.Range("A4").Value = 3
$wsFormulaCheck.Range("B4").Value = "B4"

$wsFormulaCheck.Range("A5").Value = 4
$wsFormulaCheck.Range("B5").Value = "B7"

$wsFormulaCheck.Range("A6").Value = 5
$wsFormulaCheck.Range("B6").Value = "B8"

$wsFormulaCheck.Range("A7").Value = 6
$wsFormulaCheck.Range("B7").Value = "B9"

$wsFormulaCheck.Range("A8").Select()

Write-Host "edit complete"
